$d = $word.ActiveDocument

function Get-ParaIndexForText($searchText) {
    $r = $d.Content
    $r.Find.Execute($searchText)
    return $d.Range(0, $r.Start).Paragraphs.Count + 1
}

# ---------------------------------------------------------------------
# 1) Insert a new sub-bullet (ilvl=2, numId=1) right after the paragraph
#    "Maybe look into physically interactable UI like pulling a lever."
#    holding the gamedevbeginner.com link, before "Tasks to do for this:"
# ---------------------------------------------------------------------
$idx1 = Get-ParaIndexForText("Maybe look into physically interactable UI like pulling a lever.")
$anchorPara1 = $d.Paragraphs.Item($idx1)
$anchorPara1.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item($idx1 + 1)
$newPara1.Range.ListFormat.ListIndent()
$newPara1.Range.Text = "https://gamedevbeginner.com/how-to-move-an-object-with-the-mouse-in-unity-in-2d/"

# ---------------------------------------------------------------------
# 2) Insert four new sub-bullets (ilvl=1, numId=3) right after
#    "Pullable levers and rotatable dials", before the blue
#    "Backend AI to manage ..." paragraph.
# ---------------------------------------------------------------------
$newItems = @(
    "Movable sprite with mouse",
    "Rotate object with mouse drag",
    "Fixed point at one end ",
    "Lock the pivot to be between a certain angle."
)

$idx2 = Get-ParaIndexForText("Pullable levers and rotatable dials")
$cursor = $idx2
$isFirstItem = $true
foreach ($itemText in $newItems) {
    $anchorPara = $d.Paragraphs.Item($cursor)
    $anchorPara.Range.InsertParagraphAfter()
    $cursor = $cursor + 1
    $itemPara = $d.Paragraphs.Item($cursor)
    if ($isFirstItem) {
        # Only the first new paragraph needs to be demoted one level (from
        # the inherited ilvl=0 down to ilvl=1); later paragraphs already
        # inherit ilvl=1 from the previous new paragraph.
        $itemPara.Range.ListFormat.ListIndent()
        $isFirstItem = $false
    }
    $itemPara.Range.Text = $itemText
}

# ---------------------------------------------------------------------
# 3) Merge the two "Backend AI to manage how the energy and materials "
#    + "work." runs into a single run.
# ---------------------------------------------------------------------
$mergeRng = $d.Content
$mergeRng.Find.Execute("Backend AI to manage how the energy and materials work.")
$mergeRng.Text = "IRON_TMP_PLACEHOLDER"
$mergeRng2 = $d.Content
$mergeRng2.Find.Execute("IRON_TMP_PLACEHOLDER")
$mergeRng2.Text = "Backend AI to manage how the energy and materials work."

# ---------------------------------------------------------------------
# 4) Fill in the previously-empty last bullet (ilvl=0, numId=4).
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Text = "Needs to manage all the backend stuff like production statistics like energy usage."

# ---------------------------------------------------------------------
# 5) The "ilvl=1" level of the numbering definition behind numId=3 is no
#    longer tentative now that it is actually used in the document.
# ---------------------------------------------------------------------
$idxMovable = Get-ParaIndexForText("Movable sprite with mouse")
$movablePara = $d.Paragraphs.Item($idxMovable)
$lvl2 = $movablePara.Range.ListFormat.ListTemplate.ListLevels.Item(2)
$lvl2.NumberFormat = $lvl2.NumberFormat
